$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 98: section header "Alerts"
$ws.Range("A98").Value = "Alerts"

# Rows 99-104: locator rows (ElementID, ElementPath, Method)
$data = @(
    @("alertsNav",    "//div[@class='element-list collapse show']//li[@id='item-1']"),
    @("alertsScroll", "//h1[normalize-space()='Alerts']"),
    @("alertBtn",     "//button[@id='alertButton']"),
    @("fiveMinBtn",   "//button[@id='timerAlertButton']"),
    @("confirmBoxBtn","//button[@id='confirmButton']"),
    @("promptBtn",    "//button[@id='promtButton']")
)

$row = 99
foreach ($item in $data) {
    $ws.Cells.Item($row, 1).Value = $item[0]
    $ws.Cells.Item($row, 2).Value = $item[1]
    $ws.Cells.Item($row, 3).Value = "By.xpath"
    $row++
}

# Update selection to B104
$ws.Range("B104").Select()
